$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.478.37'
$ws.Range("E2").Value = '  +2.58%  '
$ws.Range("D3").Value = '3.195.78'
$ws.Range("E3").Value = '  +1.41%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.39'
$ws.Range("E5").Value = '  +3.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.56'
$ws.Range("E6").Value = '  +4.20%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.560'
$ws.Range("E8").Value = '  +6.37%  '
$ws.Range("D9").Value = '3.195.19'
$ws.Range("E9").Value = '  +1.43%  '
$ws.Range("E10").Value = '  +1.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.89'
$ws.Range("E11").Value = '  -3.75%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.520'
$ws.Range("E12").Value = '  +3.82%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000271'
$ws.Range("E13").Value = '  +2.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '39.34'
$ws.Range("E14").Value = '  +5.92%  '
$ws.Range("D15").Value = '3.719.41'
$ws.Range("E15").Value = '  +1.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.52'
$ws.Range("E16").Value = '  +5.17%  '
$ws.Range("D17").Value = '66.474.47'
$ws.Range("E17").Value = '  +2.57%  '
$ws.Range("D18").Value = '3.197.23'
$ws.Range("E18").Value = '  +1.62%  '
$ws.Range("E19").Value = '  +0.83%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '519.41'
$ws.Range("E20").Value = '  +2.88%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.46'
$ws.Range("E21").Value = '  +3.72%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.741'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.15'
$ws.Range("E23").Value = '  +5.55%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.98'
$ws.Range("E24").Value = '  -1.27%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.99'
$ws.Range("E25").Value = '  +1.87%  '
$ws.Range("E27").Value = '  +3.71%  '
$ws.Range("E28").Value = '  +3.13%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.37'
$ws.Range("E29").Value = '  +8.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.11'
$ws.Range("E30").Value = '  +13.96%  '
$ws.Range("E31").Value = '  +5.38%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '28.40'
$ws.Range("E32").Value = '  +2.91%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.24'
$ws.Range("E33").Value = '  +3.03%  '
$ws.Range("E34").Value = '  +0.11%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.56'
$ws.Range("E35").Value = '  +1.28%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '509.66'
$ws.Range("E36").Value = '  +5.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '54.89'
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0907'
$ws.Range("E38").Value = '  +1.89%  '
$ws.Range("E40").Value = '  +9.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.94'
$ws.Range("E41").Value = '  +2.18%  '
$ws.Range("E42").Value = '  -1.04%  '
$ws.Range("D43").Value = '0.0₃0683'
$ws.Range("E43").Value = '  +16.09%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.302'
$ws.Range("E44").Value = '  +7.35%  '
$ws.Range("E45").Value = '  +1.08%  '
$ws.Range("D46").Value = '2.902.69'
$ws.Range("E46").Value = '  -2.85%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.58'
$ws.Range("E47").Value = '  +1.68%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.81'
$ws.Range("E48").Value = '  +13.67%  '
$ws.Range("E49").Value = '  +3.76%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.36'
$ws.Range("E51").Value = '  +5.48%  '
